$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new value in D6 (this adds "fio" to the shared strings table
# and extends the used range / dimension to A2:H6)
$ws.Range("D6").Value = "fio"

# Move the active selection to D6, matching the post-edit selection state
$ws.Range("D6").Select()
